# Update: add some columns for user table
# Adds two new rows ("sign_contract" and "confirm_payment_code") to the
# "user" sheet's field table, pushing the rows below them down.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("user")

# Insert a new row before row 16 (currently "registration_system"),
# and another new row before what becomes row 18 (currently "active_code").
$ws.Rows("16:16").Insert()
$ws.Rows("18:18").Insert()

# --- Row 16: sign_contract ---
$ws.Cells.Item(16, 1).Formula = "=ROW()-3"
$ws.Cells.Item(16, 2).Value = "sign_contract"
$ws.Cells.Item(16, 3).Value = "tinyint"
$ws.Cells.Item(16, 4).Value = 1

# --- Row 18: confirm_payment_code --- (string added before the H16 note below
# so the shared-string table order matches the source workbook)
$ws.Cells.Item(18, 1).Formula = "=ROW()-3"
$ws.Cells.Item(18, 2).Value = "confirm_payment_code"
$ws.Cells.Item(18, 3).Value = "varchar"
$ws.Cells.Item(18, 4).Value = 100

$ws.Cells.Item(16, 8).Value = "0: unsign, 1: signed"
$ws.Cells.Item(16, 9).Value = 0

# Expand the "Table2" structured table to include the two new rows.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A3:J23"))

$ws.Range("H16").Select()
